$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2, shifting everything else down.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new release entry.
$ws.Range("A2").Value = "Tranquillitatis"
$ws.Range("B2").Value = 45506
$ws.Range("C2").Value = "OFFN014"
$ws.Range("D2").Value = "tranquillitatis"
$ws.Range("J2").Value = "Music for moonwalkers I"
$ws.Range("K2").Value = "Offnominal"
$ws.Range("L2").Value = "https://logickal.bandcamp.com"
$ws.Range("M2").Value = "amb"

# Match the releaseDate cell's number format to the one already used for the
# catalogNumber column (numFmtId 14) instead of creating a brand-new custom
# number format.
$ws.Range("C3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("N2").Select()
